{"js": "// Merge the \"To\" + \" \" + \"Do\" heading runs into a single run reading \"ToDo\".\n// Target the heading paragraph specifically via its bookmark (the TOC entry\n// has the same visible text \"To Do\" but is a separate hyperlink range).\nconst toDoRange = context.document.getBookmarkRange(\"_Toc132530065\");\ntoDoRange.load(\"text\");\nawait context.sync();\n\nif (toDoRange.text === \"To Do\") {\n  toDoRange.insertText(\"ToDo\", \"Replace\");\n  await context.sync();\n}\n\n// Merge the \"Upon startup...\" list item runs into a single run.\nconst startupResults = context.document.body.search(\n  \"Upon startup, the main program checks to see if config file exists, if exists, load it.\",\n  { matchCase: true, matchWholeWord: false }\n);\nstartupResults.load(\"items\");\nawait context.sync();\n\nif (startupResults.items.length > 0) {\n  startupResults.items[0].insertText(\n    \"Upon startup, the main program checks to see if config file exists, if exists, load it.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Merge the \"To\" + \" \" + \"Do\" heading runs into a single run \"ToDo\". ---\n# Target the Heading2 paragraph via its bookmark so we don't touch the\n# identically-worded TOC hyperlink entry earlier in the document.\n$toDoRange = $d.Bookmarks(\"_Toc132530065\").Range\nif ($toDoRange.Text -eq \"To Do\") {\n    $toDoRange.Text = \"ToDo\"\n}\n\n# --- Merge the \"Upon startup...\" list item runs into a single run. ---\n# The resulting text is identical to the concatenation of the two existing\n# runs, so a direct Range.Text assignment to the same string is a no-op for\n# the XML (the runs would stay split). Stage the replacement through a\n# throwaway placeholder first so the real set always performs a structural\n# rewrite that merges the runs.\n$startupText = \"Upon startup, the main program checks to see if config file exists, if exists, load it.\"\n$r = $d.Content\n$found = $r.Find.Execute($startupText)\nif ($found) {\n    $r.Text = \"__TMP_PLACEHOLDER_0001__\"\n\n    $r2 = $d.Content\n    $found2 = $r2.Find.Execute(\"__TMP_PLACEHOLDER_0001__\")\n    if ($found2) {\n        $r2.Text = $startupText\n    }\n}\n"}
